$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 187; this shifts existing rows 187-236 down to 188-237
# and keeps formatting (e.g. the date-format style on column D) consistent
# with the row that was there before.
$ws.Rows("187:187").Insert()

# Populate the newly inserted row 187 with the new data record.
$ws.Range("A187").Value2 = 3
$ws.Range("B187").Value2 = "Femacal de La Calera"
$ws.Range("C187").Value2 = "Coquimbo"
$ws.Range("D187").Value2 = 44508
$ws.Range("E187").Value2 = 5
$ws.Range("F187").Value2 = 100112040
$ws.Range("G187").Value2 = "Cilantro"
$ws.Range("H187").Value2 = "Sin especificar"
$ws.Range("I187").Value2 = "Primera"
$ws.Range("J187").Value2 = 280
$ws.Range("K187").Value2 = 2000
$ws.Range("L187").Value2 = 2300
$ws.Range("M187").Value2 = 2129
$ws.Range("N187").Value2 = "`$/docena de atados (3 kilos)"
$ws.Range("O187").Value2 = "Provincia de Quillota"
$ws.Range("P187").Value2 = 710
$ws.Range("Q187").Value2 = 3
$ws.Range("R187").Value2 = "Hortaliza"
